# "done Eng Report Form" - clean up the Format_Eng sheet:
#  - clear the "jupiter ver4.1.1 / JPT-OPTISHAPE / Other" comment cells (D7:D10)
#    and let the rows shrink back to the default height
#  - clear the leftover sample data row (B22:D22)
#  - move the active selection to B23

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Format_Eng")

# Remove the boilerplate text from D7:D10 and auto-fit the rows so the
# now-unneeded extra row height (55.2) goes away.
$ws.Range("D7:D10").ClearContents()
$ws.Range("7:10").Rows.AutoFit()

# Clear the leftover example row at the bottom of the sheet.
$ws.Range("B22:D22").ClearContents()

# Leave the selection on B23, like the finished form.
$ws.Range("B23").Select()
